$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Erros - Tabela Geral" table rows (11-14) ---
# Row 11: Erros - Tabela Geral
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = "Erros - Tabela Geral"

# Row 12: Correcao Transacao Cartao
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = "Correcao Transacao Cartao"

# Row 13: Verificar Saldo no Carregamento da tela
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "Verificar Saldo no Carregamento da tela"

# Row 14: Validar Insert Transacao
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = "Validar Insert Transacao"

# --- Highlight G11 ("Erros - Tabela Geral") with a light blue fill and border ---
$g11 = $ws.Range("G11")
$g11.Interior.ThemeColor = 2
$g11.Interior.PatternColor = 16247774

$g11.Borders.Item(8).Color = 15123357
$g11.Borders.Item(9).Color = 15123357
$g11.Borders.Item(10).Color = 15123357

# --- Update active selection to E20 ---
$ws.Range("E20").Select()
